$d = $word.ActiveDocument

# 1. QUALIFICATIONS / Software line: "SolidWorks (3D modeling, Flow Simulation,
#    Finite Element Analysis), OnShape, Excel" -> "...OnShape, AutoCAD, Excel"
#    (the "(3D modeling...)" text itself is unchanged, just re-run-split upstream;
#    the actual content delta is inserting "AutoCAD, " before "Excel").
$null = $d.Content.Find.Execute(", Excel", $true, $false, $false, $false, $false, $true, 1, $false, ", AutoCAD, Excel", 2)

# 2. Bullet "Trained new AIAA members ..." paragraph: reduce space-after from
#    12pt (240 twips) to 6pt (120 twips).
$r = $d.Content.Duplicate
$null = $r.Find.Execute("Trained new AIAA members")
$r.Paragraphs.Item(1).Format.SpaceAfter = 6

# 3. "Participated in a 2-week long introduction course ..." ->
#    "Participated in a 2-week long training course ..."
$null = $d.Content.Find.Execute("Participated in a 2-week long introduction course in plasma physics and fusion energy.", $true, $false, $false, $false, $false, $true, 1, $false, "Participated in a 2-week long training course in plasma physics and fusion energy.", 2)

# 4. "Researched novel x-ray 2D dual crystal spectroscopy imaging system expected
#     to outperform previous methods ..." ->
#    "...imaging system that alleviates imaging errors from previous methods ..."
$null = $d.Content.Find.Execute("imaging system expected to outperform previous methods", $true, $false, $false, $false, $false, $true, 1, $false, "imaging system that alleviates imaging errors from previous methods", 2)

# 5. "...raytracing python package. Simulated up to 50 million photons using the
#     python package." ->
#    "...raytracing python package, which allowed for simulations of up to 50
#     million photons."
$null = $d.Content.Find.Execute("a raytracing python package. Simulated up to 50 million photons using the python package.", $true, $false, $false, $false, $false, $true, 1, $false, "a raytracing python package, which allowed for simulations of up to 50 million photons.", 2)
